$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Original layout: A time | B susceptible | C infected | D UD | E UR | F HD
#                  | G HR | H QD | I QR | J died | K immune state
#
# Insert a new column at C for "exposed" (shifts old C..K infected..immune -> D..L)
$ws.Columns("C").Insert()

# Insert a new column at L for "recovered" (shifts old immune-state column L -> M)
$ws.Columns("L").Insert()

# Insert a new trailing column at N for "sum" (inherits header style cleanly)
$ws.Columns("N").Insert()

# --- Header row ------------------------------------------------------------
$ws.Range("C1").Value = "exposed"
$ws.Range("L1").Value = "recovered"
$ws.Range("N1").Value = "sum"

# --- Data: new "exposed" (C) values, corrected "died" (K) values
# (death & recovery rates were switched), new "recovered" (L) values,
# and new "sum" (N) values -------------------------------------------------
$ws.Range("C2").Value = 0.7013472200632157
$ws.Range("K2").Value = 0.00002901704751541531
$ws.Range("L2").Value = 0.1180921291258615
$ws.Range("N2").Value = 1.701347220063216
$ws.Range("C3").Value = 0.6995826969419586
$ws.Range("K3").Value = 0.00002902528615984391
$ws.Range("L3").Value = 0.0004112392469746046
$ws.Range("N3").Value = 1.493838456037322
$ws.Range("C4").Value = 0.7152828071409932
$ws.Range("K4").Value = 0.00003271186085128183
$ws.Range("L4").Value = 0.0005266132234678202
$ws.Range("N4").Value = 1.480201601711222
$ws.Range("C5").Value = 0.7397431426785109
$ws.Range("K5").Value = 0.0000411368965342757
$ws.Range("L5").Value = 0.0006336332576286177
$ws.Range("N5").Value = 1.462682770515222
$ws.Range("C6").Value = 0.7663179722163316
$ws.Range("K6").Value = 0.00005515205236754834
$ws.Range("L6").Value = 0.0007290662323813206
$ws.Range("N6").Value = 1.441573086559354
$ws.Range("C7").Value = 0.7903091162295177
$ws.Range("K7").Value = 0.00007547587986205919
$ws.Range("L7").Value = 0.0008167274476885717
$ws.Range("N7").Value = 1.417265341776575
$ws.Range("C8").Value = 0.8086746965839954
$ws.Range("K8").Value = 0.0001027353591913873
$ws.Range("L8").Value = 0.0008983480896995168
$ws.Range("N8").Value = 1.39002203282684
$ws.Range("C9").Value = 0.8197141702734354
$ws.Range("K9").Value = 0.0001374843512411309
$ws.Range("L9").Value = 0.0009741542340481884
$ws.Range("N9").Value = 1.360046111185132
$ws.Range("C10").Value = 0.8227765460717776
$ws.Range("K10").Value = 0.0001802076620112761
$ws.Range("L10").Value = 0.001043526145171184
$ws.Range("N10").Value = 1.327532054330388
$ws.Range("C11").Value = 0.8179965127169085
$ws.Range("K11").Value = 0.0002313172527893831
$ws.Range("L11").Value = 0.00110548214486326
$ws.Range("N11").Value = 1.292694668459822
$ws.Range("C12").Value = 0.8060560036011692
$ws.Range("K12").Value = 0.0002911452568505292
$ws.Range("L12").Value = 0.00115901252268354
$ws.Range("N12").Value = 1.255781536204499
$ws.Range("C13").Value = 0.7879731321595217
$ws.Range("K13").Value = 0.000359936953357352
$ws.Range("L13").Value = 0.00120329317772213
$ws.Range("N13").Value = 1.217074163201525
$ws.Range("C14").Value = 0.7649244771341219
$ws.Range("K14").Value = 0.0004378456565185131
$ws.Range("L14").Value = 0.00123780377954357
$ws.Range("N14").Value = 1.176881877073638
$ws.Range("C15").Value = 0.7381065413739042
$ws.Range("K15").Value = 0.000524930546697753
$ws.Range("L15").Value = 0.001262372034323321
$ws.Range("N15").Value = 1.135531698377883
$ws.Range("C16").Value = 0.7086385473477969
$ws.Range("K16").Value = 0.0006211577632845155
$ws.Range("L16").Value = 0.001277163991553548
$ws.Range("N16").Value = 1.093356678052296
$ws.Range("C17").Value = 0.6775039059920022
$ws.Range("K17").Value = 0.0007264045753143201
$ws.Range("L17").Value = 0.00128263888858664
$ws.Range("N17").Value = 1.050684526509436
$ws.Range("C18").Value = 0.6455237990670438
$ws.Range("K18").Value = 0.0008404661261269205
$ws.Range("L18").Value = 0.001279484839469663
$ws.Range("N18").Value = 1.007827737699688
$ws.Range("C19").Value = 0.6133543796971874
$ws.Range("K19").Value = 0.0009630640874043328
$ws.Range("L19").Value = 0.001268548529928728
$ws.Range("N19").Value = 0.9650758621771519
$ws.Range("C20").Value = 0.5814991050248772
$ws.Range("K20").Value = 0.001093856522311167
$ws.Range("L20").Value = 0.001250768369074226
$ws.Range("N20").Value = 0.9226901387424405
$ws.Range("C21").Value = 0.5503290470256248
$ws.Range("K21").Value = 0.001232448309494453
$ws.Range("L21").Value = 0.001227116875375494
$ws.Range("N21").Value = 0.8809003741591631
$ws.Range("C22").Value = 0.5201059290371914
$ws.Range("K22").Value = 0.001378401583164016
$ws.Range("L22").Value = 0.001198554951317933
$ws.Range("N22").Value = 0.8399037639447255
